# Apply the "Added rest assured framework code" edit:
#  - rename the existing sheet to cxCreationValidKey and trim its test data
#    down to 3 engineers (rows 2-4), dropping engineers 4-9
#  - add a new sheet cxCreationInvalidKey with the first two engineers but
#    a different description column, and make it the active sheet
#  - re-wire the mailto hyperlinks on both sheets to match the new data

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet 1: createCustomerWithValidKey -> cxCreationValidKey ---
$ws1.Name = "cxCreationValidKey"

# Hyperlinks.Delete() on this host removes every hyperlink on the sheet
# (scoping is ignored), so clear them up front and re-add the ones we
# still need once the stale rows are gone.
$ws1.Hyperlinks.Delete()

# Drop rows 5-10 (testEngineer4..testEngineer9), keeping header + 3 rows.
$ws1.Range("A5:A10").EntireRow.Delete()

$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:test1@mailinator.com")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:test2@mailinator.com")
$ws1.Hyperlinks.Add($ws1.Range("B4"), "mailto:test3@mailinator.com")
$ws1.Range("B2:B4").Style = "Hyperlink"

$ws1.Range("C14").Select()

# --- Sheet 2: new cxCreationInvalidKey sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "cxCreationInvalidKey"

$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "email"
$ws2.Range("C1").Value = "description"

$ws2.Range("A2").Value = "testEngineer1"
$ws2.Range("B2").Value = "test1@mailinator.com"
$ws2.Range("C2").Value = "2ndSheetDescriptionTestEngineer1"

$ws2.Range("A3").Value = "testEngineer2"
$ws2.Range("B3").Value = "test2@mailinator.com"
$ws2.Range("C3").Value = "2ndSheetDescriptionTestEngineer2"

$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:test1@mailinator.com")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:test2@mailinator.com")
$ws2.Range("B2:B3").Style = "Hyperlink"

$ws2.Range("D10").Select()

# cxCreationInvalidKey ends up the active/selected tab, matching the diff's
# activeTab="1" / tabSelected="1" on the second sheet.
$ws2.Activate()
